# Adds a new "culture_collection" field/column to the MIGS.ba plant-associated
# template, inserted (alphabetically) right before the existing "depth" column.
#
# Concretely, on row 15 (the header row) a new column is inserted immediately
# before column Y ("depth"), pushing "depth", "elev", ... all one column to
# the right (Y->Z, Z->AA, ... CH->CI). The new column Y gets the header
# "culture_collection" plus its own cell comment; every comment that lived in
# columns Y..CH has to move one column to the right along with its cell's
# content (EntireColumn.Insert shifts cell values/formatting automatically,
# but not the legacy cell-comment objects, so those are re-applied by hand
# below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColName([int]$n) {
    $s = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $s = [char](65 + $rem) + $s
        $n = [int](($n - 1) / 26)
    }
    return $s
}

$firstCol = 25   # column Y (1-based): the new column is inserted here
$lastCol  = 86   # column CH (1-based): last existing column on row 15

# --- Step 1: capture all existing comment texts for columns Y..CH (row 15) ---
# before the insert shifts their underlying cell content out from under them.
$texts = @{}
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $ref = (ColName $col) + "15"
    $cm = $ws.Range($ref).Comment
    if ($cm -ne $null) {
        $texts[$col] = $cm.Text()
    } else {
        $texts[$col] = $null
    }
}

# --- Step 2: insert the new column before Y; cell values/styles shift right ---
$ws.Range((ColName $firstCol) + "1").EntireColumn.Insert()

# --- Step 3: populate the new header cell ---
$ws.Range((ColName $firstCol) + "15").Value = "culture_collection"

# --- Step 4: give the new column its own comment ---
$cultureCollectionComment = "Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier"
$newCell = $ws.Range((ColName $firstCol) + "15")
if ($newCell.Comment -ne $null) {
    $newCell.Comment.Text($cultureCollectionComment) | Out-Null
} else {
    $newCell.AddComment($cultureCollectionComment) | Out-Null
}

# --- Step 5: re-home the captured comments one column to the right ---
# (iterate from the right so we never clobber a comment before reading it)
for ($col = $lastCol + 1; $col -ge $firstCol + 1; $col--) {
    $ref = (ColName $col) + "15"
    $cell = $ws.Range($ref)
    $srcText = $texts[$col - 1]
    if ($srcText -ne $null) {
        if ($cell.Comment -ne $null) {
            $cell.Comment.Text($srcText) | Out-Null
        } else {
            $cell.AddComment($srcText) | Out-Null
        }
    } else {
        if ($cell.Comment -ne $null) {
            $cell.Comment.Delete() | Out-Null
        }
    }
}
